$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (preserve rich-text runs) ---
$hdr = $ws.Range("A8")
$hdr.Characters(21, 2).Text = "37"

$rpt = $ws.Range("C9")
$rpt.Characters(27, 8).Text = "9/9/2024"
$rpt.Characters(46, 8).Text = "9/15/2024"

# --- Data table updates (rows 15-31), new weekly crime figures ---
# Donor cells (stable, style-only copies used to fix up number format / string-type
# after a cell switches between numeric and text representation).
$donorText = "C15"   # style s=14 (text)
$donorNum1 = "J14"   # style s=15 (integer count)
$donorNum2 = "K14"   # style s=16 (percentage)

$ws.Range("N15").Value = -62.962962962963
$ws.Range("C16").Value = 1
$ws.Range("E16").Value = -66.666666666666
$ws.Range("F16").Value = 8
$ws.Range("H16").Value = -38.461538461538
$ws.Range("I16").Value = 104
$ws.Range("J16").Value = 89
$ws.Range("K16").Value = 16.853932584269
$ws.Range("L16").Value = -14.754098360655
$ws.Range("M16").Value = -26.760563380281
$ws.Range("N16").Value = -83.492063492063
$ws.Range("C17").Value = 1
$ws.Range("D17").Value = 2
$ws.Range("E17").Value = -50
$ws.Range("F17").Value = 9
$ws.Range("G17").Value = 13
$ws.Range("H17").Value = -30.769230769230
$ws.Range("I17").Value = 105
$ws.Range("J17").Value = 135
$ws.Range("K17").Value = -22.222222222222
$ws.Range("L17").Value = -20.454545454545
$ws.Range("M17").Value = 40
$ws.Range("N17").Value = -70.505617977528
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = 50
$ws.Range("F18").Value = 9
$ws.Range("G18").Value = 8
$ws.Range("H18").Value = 12.5
$ws.Range("I18").Value = 81
$ws.Range("J18").Value = 125
$ws.Range("K18").Value = -35.2
$ws.Range("L18").Value = -39.097744360902
$ws.Range("M18").Value = -10.989010989011
$ws.Range("N18").Value = -89.271523178808
$ws.Range("C19").Value = 4
$ws.Range("D19").Value = 9
$ws.Range("E19").Value = -55.555555555555
$ws.Range("F19").Value = 45
$ws.Range("G19").Value = 29
$ws.Range("H19").Value = 55.172413793103
$ws.Range("I19").Value = 350
$ws.Range("J19").Value = 362
$ws.Range("K19").Value = -3.314917127071
$ws.Range("L19").Value = -6.417112299465
$ws.Range("M19").Value = -4.891304347826
$ws.Range("N19").Value = -54.007884362680
$ws.Range("C20").Value = 3
$ws.Range($donorNum1).Copy()
$ws.Range("C20").PasteSpecial(-4122)
$ws.Range("D20").Value = 4
$ws.Range("E20").Value = -25
$ws.Range("F20").Value = 9
$ws.Range("G20").Value = 8
$ws.Range("H20").Value = 12.5
$ws.Range("I20").Value = 48
$ws.Range("J20").Value = 75
$ws.Range("K20").Value = -36
$ws.Range("L20").Value = -26.153846153846
$ws.Range("M20").Value = 77.777777777777
$ws.Range("N20").Value = -93.172119487909
$ws.Range("C21").Value = 12
$ws.Range("D21").Value = 20
$ws.Range("E21").Value = -40
$ws.Range("F21").Value = 80
$ws.Range("G21").Value = 72
$ws.Range("H21").Value = 11.111111111111
$ws.Range("I21").Value = 698
$ws.Range("J21").Value = 798
$ws.Range("K21").Value = -12.531328320802
$ws.Range("L21").Value = -16.105769230769
$ws.Range("M21").Value = -2.103786816269
$ws.Range("N21").Value = -78.516466605109
$ws.Range("C22").Value = "'0"
$ws.Range($donorText).Copy()
$ws.Range("C22").PasteSpecial(-4122)
$ws.Range("D22").Value = "'0"
$ws.Range($donorText).Copy()
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("E22").Value = "'***.*"
$ws.Range($donorText).Copy()
$ws.Range("E22").PasteSpecial(-4122)
$ws.Range("G22").Value = 1
$ws.Range("H22").Value = 300
$ws.Range("L22").Value = 15.789473684210
$ws.Range("M22").Value = 22.222222222222
$ws.Range("C23").Value = 1
$ws.Range("D23").Value = 2
$ws.Range("E23").Value = -50
$ws.Range("G23").Value = 9
$ws.Range("H23").Value = -33.333333333333
$ws.Range("I23").Value = 69
$ws.Range("J23").Value = 82
$ws.Range("K23").Value = -15.853658536585
$ws.Range("L23").Value = -24.175824175824
$ws.Range("M23").Value = 21.052631578947
$ws.Range("C24").Value = 21
$ws.Range("D24").Value = 26
$ws.Range("E24").Value = -19.230769230769
$ws.Range("F24").Value = 83
$ws.Range("G24").Value = 102
$ws.Range("H24").Value = -18.627450980392
$ws.Range("I24").Value = 817
$ws.Range("J24").Value = 1147
$ws.Range("K24").Value = -28.770706190061
$ws.Range("L24").Value = -37.776085300837
$ws.Range("M24").Value = 10.554803788903
$ws.Range("C25").Value = 8
$ws.Range("D25").Value = 18
$ws.Range("E25").Value = -55.555555555555
$ws.Range("F25").Value = 41
$ws.Range("G25").Value = 77
$ws.Range("H25").Value = -46.753246753246
$ws.Range("I25").Value = 442
$ws.Range("J25").Value = 789
$ws.Range("K25").Value = -43.979721166033
$ws.Range("L25").Value = -47.754137115839
$ws.Range("C26").Value = 11
$ws.Range("E26").Value = 57.142857142857
$ws.Range("F26").Value = 29
$ws.Range("G26").Value = 19
$ws.Range("H26").Value = 52.631578947368
$ws.Range("I26").Value = 231
$ws.Range("J26").Value = 212
$ws.Range("K26").Value = 8.962264150943
$ws.Range("L26").Value = 8.450704225352
$ws.Range("M26").Value = -3.75
$ws.Range("D27").Value = "'0"
$ws.Range($donorText).Copy()
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("E27").Value = "'***.*"
$ws.Range($donorText).Copy()
$ws.Range("E27").PasteSpecial(-4122)
$ws.Range("C28").Value = "'0"
$ws.Range($donorText).Copy()
$ws.Range("C28").PasteSpecial(-4122)
$ws.Range("D28").Value = 1
$ws.Range($donorNum1).Copy()
$ws.Range("D28").PasteSpecial(-4122)
$ws.Range("E28").Value = -100
$ws.Range($donorNum2).Copy()
$ws.Range("E28").PasteSpecial(-4122)
$ws.Range("F28").Value = 3
$ws.Range("H28").Value = 50
$ws.Range("J28").Value = 31
$ws.Range("K28").Value = -12.903225806451
$ws.Range("D31").Value = 1
$ws.Range($donorNum1).Copy()
$ws.Range("D31").PasteSpecial(-4122)
$ws.Range("E31").Value = -100
$ws.Range($donorNum2).Copy()
$ws.Range("E31").PasteSpecial(-4122)
$ws.Range("G31").Value = 1
$ws.Range($donorNum1).Copy()
$ws.Range("G31").PasteSpecial(-4122)
$ws.Range("H31").Value = -100
$ws.Range($donorNum2).Copy()
$ws.Range("H31").PasteSpecial(-4122)
$ws.Range("J31").Value = 11
$ws.Range("K31").Value = 54.545454545454
